$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '90.606.15'
$ws.Cells.Item(2, 5).Value = '  +1.10%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.190.51'
$ws.Cells.Item(3, 5).Value = '  +4.14%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.25%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '239.70'
$ws.Cells.Item(5, 5).Value = '  +1.29%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '618.81'
$ws.Cells.Item(6, 5).Value = '  +0.39%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +5.93%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.373'
$ws.Cells.Item(8, 5).Value = '  +2.25%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.12%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '3.173.88'
$ws.Cells.Item(10, 5).Value = '  +3.60%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +5.51%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +1.08%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +0.25%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '35.23'
$ws.Cells.Item(14, 5).Value = '  +1.81%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '5.57'
$ws.Cells.Item(15, 5).Value = '  +3.92%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '90.714.23'
$ws.Cells.Item(16, 5).Value = '  +1.51%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.751.28'
$ws.Cells.Item(17, 5).Value = '  +3.13%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.146.49'
$ws.Cells.Item(18, 5).Value = '  +2.81%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.72'
$ws.Cells.Item(19, 5).Value = '  -1.66%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '15.15'
$ws.Cells.Item(20, 5).Value = '  +10.65%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.02'
$ws.Cells.Item(21, 5).Value = '  +11.69%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '451.28'
$ws.Cells.Item(22, 5).Value = '  +5.05%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -3.86%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.18'
$ws.Cells.Item(24, 5).Value = '  +5.90%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '5.75'
$ws.Cells.Item(25, 5).Value = '  +3.74%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '11.98'
$ws.Cells.Item(26, 5).Value = '  +2.89%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '82.64'
$ws.Cells.Item(27, 5).Value = '  +1.43%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).Value = '  -0.24%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.141'
$ws.Cells.Item(30, 5).Value = '  +58.05%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +19.50%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.172'
$ws.Cells.Item(32, 5).Value = '  +8.53%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '9.35'
$ws.Cells.Item(33, 5).Value = '  +4.78%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.172'
$ws.Cells.Item(34, 5).Value = '  +14.87%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -5.88%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '26.71'
$ws.Cells.Item(36, 5).Value = '  +4.80%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '7.68'
$ws.Cells.Item(37, 5).Value = '  +8.37%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Bittensor'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '510.27'
$ws.Cells.Item(38, 5).Value = '  +4.36%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'PancakeSwap'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.98'
$ws.Cells.Item(39, 5).Value = '  +5.79%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.35'
$ws.Cells.Item(40, 5).Value = '  +7.84%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.451'
$ws.Cells.Item(41, 5).Value = '  +13.35%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'MantraDAO'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.83'
$ws.Cells.Item(42, 5).Value = '  -8.60%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.45'
$ws.Cells.Item(43, 5).Value = '  -4.32%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '22.06'
$ws.Cells.Item(44, 5).Value = '  -0.05%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.729'
$ws.Cells.Item(46, 5).Value = '  +8.93%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.94'
$ws.Cells.Item(47, 5).Value = '  +4.71%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '156.61'
$ws.Cells.Item(48, 5).Value = '  +0.10%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +7.25%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '4.47'
$ws.Cells.Item(50, 5).Value = '  +3.92%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '43.96'
$ws.Cells.Item(51, 5).Value = '  -1.10%  '
